$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# Add two new library entries (rows 68-69) plus a trailing blank row
# (row 70), mirroring the existing "Spring-boot" block (rows 60-67).
# We clone row 67's formatting (font/alignment/row-height -> style 15)
# via Copy + Insert so the new rows keep the same cell style index
# instead of Excel fabricating an equivalent-but-different style.
# ---------------------------------------------------------------------

# Row 68 - new row cloned from row 67
$ws.Range("A67:C67").Copy()
$ws.Range("A68:C68").Insert()
$ws.Rows.Item(68).RowHeight = 33

# Row 69 - new row cloned from row 67
$ws.Range("A67:C67").Copy()
$ws.Range("A69:C69").Insert()
$ws.Rows.Item(69).RowHeight = 33

# Row 70 - trailing (mostly empty) row, only column A carries the style
$ws.Range("A67").Copy()
$ws.Range("A70").Insert()
$ws.Rows.Item(70).RowHeight = 33
$ws.Range("B70:C70").Clear()
$ws.Range("A70").ClearContents()

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Set the cell text in the same order the original authoring tool would
# have introduced new shared strings, so that new <si> entries land at
# indices 223-226 in the same order as the reference edit.
# ---------------------------------------------------------------------

# Row 68, col A keeps "Spring-boot" (already cloned from row 67)
$ws.Range("A68").Value = "Spring-boot"
# -> new shared string #223
$ws.Range("B68").Value = "Class CommandLineRunner and @Order"

# Row 69
$ws.Range("A69").Value = "Spring"
# -> new shared string #224
$ws.Range("B69").Value = "Use java-config with @Bean"
# -> new shared string #225
$ws.Range("C69").Value = @'
The xml-config below is equivalence to java-config:
-----Java based config ---
@Configuration
public class HelloWorldConfig {
   @Bean 
   public HelloWorld helloWorld(){
      return new HelloWorld();
   }
}
-----XML based config ----
<beans>
   <bean id = "helloWorld" class = "com.tutorialspoint.HelloWorld" />
</beans>
'@

# -> new shared string #226 (set last so it lands after #225)
$ws.Range("C68").Value = @'
Interface used to indicate that a bean should run when it is contained within a SpringApplication. Multiple CommandLineRunner beans can be defined within the same application context and can be ordered using the Ordered interface or @Order annotation.
----example----
@Order(value=3)
@Component
class ApplicationStartupRunnerOne implements CommandLineRunner {
    protected final Log logger = LogFactory.getLog(getClass());
    @Override
    public void run(String... args) throws Exception {
        logger.info("ApplicationStartupRunnerOne run method Started !!");
    }
}
@Order(value=2)
@Component
class ApplicationStartupRunnerTwo implements CommandLineRunner {
    protected final Log logger = LogFactory.getLog(getClass());
    @Override
    public void run(String... args) throws Exception {
        logger.info("ApplicationStartupRunnerTwo run method Started !!");
    }
}
'@

# ---------------------------------------------------------------------
# Restore the fixed row height (setting .Value triggers Excel's
# autofit-on-wrap, which otherwise grows the row for the long text).
# ---------------------------------------------------------------------
$ws.Rows.Item(68).RowHeight = 33
$ws.Rows.Item(69).RowHeight = 33
$ws.Rows.Item(70).RowHeight = 33

# ---------------------------------------------------------------------
# Bring the view in line with the new bottom of the sheet.
# ---------------------------------------------------------------------
$ws.Range("C68").Select()
$excel.ActiveWindow.ScrollRow = 64
